$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 177
$ws.Range("F3").Value = 396
$ws.Range("F4").Value = 1105
$ws.Range("F5").Value = 34
$ws.Range("F7").Value = 10
$ws.Range("F8").Value = 1051
$ws.Range("F9").Value = 522
$ws.Range("F10").Value = 309
$ws.Range("F11").Value = 407
$ws.Range("F12").Value = 30
$ws.Range("F13").Value = 297
$ws.Range("F14").Value = 343
$ws.Range("F15").Value = 16
$ws.Range("F16").Value = 58
$ws.Range("F17").Value = 374
$ws.Range("F18").Value = 428
$ws.Range("F19").Value = 5456
$ws.Range("F21").Value = 1528
$ws.Range("F22").Value = 347
$ws.Range("F23").Value = 4584
$ws.Range("F24").Value = 4585
$ws.Range("F25").Value = 113
$ws.Range("F26").Value = 81
$ws.Range("F27").Value = 1455
$ws.Range("F29").Value = 20
$ws.Range("F30").Value = 634
$ws.Range("F31").Value = 10

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 99
$ws.Range("F10").Value = 14
$ws.Range("F16").Value = 52

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2110

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2110
$ws.Range("F5").Value = 177
$ws.Range("F6").Value = 396
$ws.Range("F7").Value = 1105
$ws.Range("F8").Value = 34
$ws.Range("F10").Value = 10
$ws.Range("F11").Value = 1051
$ws.Range("F12").Value = 522
$ws.Range("F13").Value = 309
$ws.Range("F14").Value = 407
$ws.Range("F15").Value = 30
$ws.Range("F16").Value = 297
$ws.Range("F17").Value = 343
$ws.Range("F18").Value = 16
$ws.Range("F19").Value = 58
$ws.Range("F23").Value = 374
$ws.Range("F24").Value = 428
$ws.Range("F25").Value = 5455
$ws.Range("F27").Value = 1528
$ws.Range("F30").Value = 347
$ws.Range("F32").Value = 4585
$ws.Range("F33").Value = 4585
$ws.Range("F34").Value = 113
$ws.Range("F35").Value = 81
$ws.Range("F36").Value = 1455
$ws.Range("F38").Value = 20
$ws.Range("F39").Value = 634
$ws.Range("F40").Value = 10
$ws.Range("F41").Value = 14
$ws.Range("F49").Value = 52
